$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "MDE" source-note block that currently sits at A20:A21 ---
# (it gets relocated further down the sheet, below the new table)
$ws.Range("A20").Clear()
$ws.Range("A21").Clear()

# --- New table header (row 17): Number of employees / Assets / Turnover ---
# ("title" cell style = bold Calibri 11, no other formatting)
$ws.Range("B17").Value = "Number of employees"
$ws.Range("C17").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D17").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B17:D17").Font.Bold = $true

# --- Row 18: Micro ---
$ws.Range("A18").Value = "Micro"
$ws.Range("B18").Value = "<10"
$ws.Range("D18").Value = "< USD 250,000"

# --- Row 19: Small ---
$ws.Range("A19").Value = "Small"
$ws.Range("B19").Value = "10-100"
$ws.Range("D19").Value = "USD 250,000 - 3 Millionlion"

# --- Row 20: Medium ---
$ws.Range("A20").Value = "Medium"
$ws.Range("B20").Value = "100-200"
$ws.Range("D20").Value = "USD 3 Millionlion - 10 Millionlion"

# --- Row 21: Large ---
$ws.Range("A21").Value = "Large"
$ws.Range("B21").Value = ">200"
$ws.Range("D21").Value = "> 10 Millionlion"

# --- Relocated source-note block, now at A26:A27 ---
# ("title" = bold Calibri 11; "source" = italic Calibri 11)
$ws.Range("A26").Value = "MDE"
$ws.Range("A26").Font.Bold = $true
$ws.Range("A27").Value = 'Ministeriio da Economia (MDE), "Síntese do Programa de Desenvolvimento das MPME`s", 2012, p. 6. Avaialable at http://www.minec.gov.ao/VerPublicacao.aspx?id=820'
$ws.Range("A27").Font.Italic = $true
